# Update "Periodo Mora" column (E) on Hoja1 so the period values are
# listed in descending order (newest period first) instead of ascending.
# Previously the periods ran 1705 -> 2003 top to bottom (rows 16-50);
# now they run 2003 -> 1705 top to bottom, i.e. the previous-periods list
# is rebuilt/reordered ("Elimna EC anteriores y se agregan nuevos, se
# modifica base de datos").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2003"
$ws.Range("E17").Value = "2002"
$ws.Range("E18").Value = "2001"
$ws.Range("E19").Value = "1912"
$ws.Range("E20").Value = "1911"
$ws.Range("E21").Value = "1910"
$ws.Range("E22").Value = "1909"
$ws.Range("E23").Value = "1908"
$ws.Range("E24").Value = "1907"
$ws.Range("E25").Value = "1906"
$ws.Range("E26").Value = "1905"
$ws.Range("E27").Value = "1904"
$ws.Range("E28").Value = "1903"
$ws.Range("E29").Value = "1902"
$ws.Range("E30").Value = "1901"
$ws.Range("E31").Value = "1812"
$ws.Range("E32").Value = "1811"
$ws.Range("E33").Value = "1810"
$ws.Range("E34").Value = "1809"
$ws.Range("E35").Value = "1808"
$ws.Range("E36").Value = "1807"
$ws.Range("E37").Value = "1806"
$ws.Range("E38").Value = "1805"
$ws.Range("E39").Value = "1804"
$ws.Range("E40").Value = "1803"
$ws.Range("E41").Value = "1802"
$ws.Range("E42").Value = "1801"
$ws.Range("E43").Value = "1712"
$ws.Range("E44").Value = "1711"
$ws.Range("E45").Value = "1710"
$ws.Range("E46").Value = "1709"
$ws.Range("E47").Value = "1708"
$ws.Range("E48").Value = "1707"
$ws.Range("E49").Value = "1706"
$ws.Range("E50").Value = "1705"

$wb.Save()
